$wb = $excel.ActiveWorkbook

# --- Leaderboard sheet: append two new result rows ---
$ws = $wb.Worksheets.Item("Leaderboard")

$ws.Range("A5").Value = "Ben"
$ws.Range("B5").Value = "Trophy Course"
$ws.Range("C5").Value = "2S"
$ws.Range("D5").Value = (Get-Date -Year 2025 -Month 11 -Day 23 -Hour 15 -Minute 39 -Second 0)
$ws.Range("E5").Value = 20.010000000000002

$ws.Range("A6").Value = "Ben"
$ws.Range("B6").Value = "Trophy Course"
$ws.Range("C6").Value = "2S"
$ws.Range("D6").Value = (Get-Date -Year 2025 -Month 11 -Day 23 -Hour 14 -Minute 39 -Second 0)
$ws.Range("E6").Value = 19.010000000000002

# --- Participant sheet: move the active selection ---
$wsParticipant = $wb.Worksheets.Item("Participant")
$wsParticipant.Activate()
$wsParticipant.Range("A3").Select()

# --- Restore Leaderboard as the active sheet with its new selection ---
$ws.Activate()
$ws.Range("B8").Select()
